# "update links for materialize"
#
# 1. RESUMES sheet: Joana (row 13) gets her resume PDF hyperlinked in column A
#    (previously blank).
# 2. CSS sheet: two new rows are appended (materialize.css / materialize.js)
#    each with a hyperlink in column A pointing at the iGEM wiki raw-file
#    endpoint, and the friendly filename in column B.
# 3. The active sheet/selection moves from RESUMES to CSS (selection C12),
#    and RESUMES' own remembered selection becomes A13.

$wb = $excel.ActiveWorkbook

# --- RESUMES: hyperlink Joana's resume -------------------------------------
$wsResumes = $wb.Worksheets.Item("RESUMES")

$wsResumes.Hyperlinks.Add(
    $wsResumes.Range("A13"),
    "https://2019.igem.org/wiki/images/3/3b/T--CMUQ--wikijoanaresume.pdf",
    "",
    "",
    "https://2019.igem.org/wiki/images/3/3b/T--CMUQ--wikijoanaresume.pdf"
) | Out-Null

# Hyperlinks.Add stamps its own default "Hyperlink" style - restore the
# sheet's usual linked-cell look (same as A14, the row right below) instead.
$wsResumes.Range("A14").Copy()
$wsResumes.Range("A13").PasteSpecial(-4122) | Out-Null

# --- CSS: add materialize.css / materialize.js rows -------------------------
$wsCss = $wb.Worksheets.Item("CSS")

$wsCss.Range("B7").Value = "materialize.css"
$wsCss.Hyperlinks.Add(
    $wsCss.Range("A7"),
    "https://2019.igem.org/wiki/index.php?title=Team:CMUQ/materialize-css&action=raw&ctype=text/css",
    "",
    "",
    "https://2019.igem.org/wiki/index.php?title=Team:CMUQ/materialize-css&action=raw&ctype=text/css"
) | Out-Null

$wsCss.Range("B8").Value = "materialize.js"
$wsCss.Hyperlinks.Add(
    $wsCss.Range("A8"),
    "https://2019.igem.org/wiki/index.php?title=Team:CMUQ/materialize-js&action=raw&ctype=text/javascript",
    "",
    "",
    "https://2019.igem.org/wiki/index.php?title=Team:CMUQ/materialize-js&action=raw&ctype=text/javascript"
) | Out-Null

# Re-apply the same row format as the rest of the link table (row 2) so the
# new rows match s="11"/s="5" instead of the default hyperlink style.
$wsCss.Range("A2:B2").Copy()
$wsCss.Range("A7:B7").PasteSpecial(-4122) | Out-Null
$wsCss.Range("A2:B2").Copy()
$wsCss.Range("A8:B8").PasteSpecial(-4122) | Out-Null

# --- Selection / active sheet bookkeeping -----------------------------------
# Record RESUMES' new remembered selection before switching away from it.
$wsResumes.Range("A13").Select() | Out-Null

# CSS becomes the active/visible sheet, selection parked on C12.
$wsCss.Activate() | Out-Null
$wsCss.Range("C12").Select() | Out-Null
